$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-27 Sunday" "2023-08-28 Monday"

Replace-Text "28×67=" "32×70="
Replace-Text "81×54=" "72×42="
Replace-Text "74×46=" "40×55="
Replace-Text "88×57=" "15×68="
Replace-Text "58×88=" "41×46="
Replace-Text "36×59=" "50×14="
Replace-Text "53×21=" "49×48="
Replace-Text "12×37=" "22×37="
Replace-Text "81×63=" "64×26="
Replace-Text "63×53=" "17×78="
Replace-Text "11×73=" "85×57="
Replace-Text "53×36=" "14×88="
Replace-Text "70×60=" "38×26="
Replace-Text "30×47=" "60×65="
Replace-Text "81×62=" "54×11="
Replace-Text "91×65=" "21×78="
Replace-Text "77×26=" "48×20="
Replace-Text "53×26=" "15×96="
Replace-Text "68×93=" "15×57="
Replace-Text "57×42=" "65×65="
Replace-Text "66×79=" "51×80="
Replace-Text "79×31=" "22×83="
Replace-Text "72×55=" "13×24="
Replace-Text "16×12=" "64×87="
Replace-Text "68×28=" "27×11="
